$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("line_imp")

$ws.Range("D1").Value = "shunt_r"
$ws.Range("E1").Value = "shunt_x"

$ws.Range("D2").Value = 0
$ws.Range("E2").Value = 0
$ws.Range("D3").Value = 0
$ws.Range("E3").Value = 0
$ws.Range("D4").Value = 0
$ws.Range("E4").Value = 0

$ws.Range("E5").Select()
